# Set workbook to manual calculation mode (calcPr calcMode="manual")
$excel.Calculation = -4135

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row (row 33) for a Mac-Address / Document Types entry
$ws.Range("A33").Value = 10002
$ws.Range("B33").Value = 110032
$ws.Range("C33").Value = 10032
$ws.Range("D33").Value = "eng"
$ws.Range("E33").Value = $true
$ws.Range("F33").Value = "superadmin"
$ws.Range("G33").Value = "now()"
$ws.Range("H33").Value = "now()"

# Update the active selection on the sheet
$ws.Range("E31").Select() | Out-Null
